# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, thin border, centered alignment)
# from the last existing header cell (AC1) onto the three new header cells
# so they match the rest of the header row, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-48) with the team's
# overall record for the season.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 106
    $ws.Cells.Item($r, 31).Value = 56
    $ws.Cells.Item($r, 32).Value = 0
}
